$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Investor *" column header to "Stakeholder *"
$ws.Range("A1").Value = "Stakeholder *"

# Remove the stray formatted cell far outside the real data (row 36 / column O)
# so the sheet's used range shrinks back down to the actual data (A1:M7).
$ws.Range("O36").Clear()

# Reset the active selection to A2
$ws.Range("A2").Select()
